# Weekly fruit/vegetable price update: insert a new weekly record as the
# new first historical row for this item (row 121), shifting the existing
# rows 121-126 down to 122-127.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 121; existing rows 121-126 shift down to 122-127.
$ws.Rows.Item(121).Insert()

# Populate the newly inserted row 121 with the new weekly data point.
$ws.Cells.Item(121, 1).Value = 1
$ws.Cells.Item(121, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(121, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(121, 4).Value = 45021
$ws.Cells.Item(121, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(121, 5).Value = 15
$ws.Cells.Item(121, 6).Value = 100112038
$ws.Cells.Item(121, 7).Value = "Cebollín baby"
$ws.Cells.Item(121, 8).Value = "Sin especificar"
$ws.Cells.Item(121, 9).Value = "Primera"
$ws.Cells.Item(121, 10).Value = 300
$ws.Cells.Item(121, 11).Value = 2000
$ws.Cells.Item(121, 12).Value = 2500
$ws.Cells.Item(121, 13).Value = 2250
$ws.Cells.Item(121, 14).Value = "`$/paquete 1,5 a 2 kilos"
$ws.Cells.Item(121, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(121, 16).Value = 1125
$ws.Cells.Item(121, 17).Value = 2
$ws.Cells.Item(121, 18).Value = "Hortaliza"
